# Journal policies.docx - F22 update
# Applies the edits described by the commit "update of policy docs for F22":
#   1. Insert a new reminder paragraph after "Journals will be submitted..."
#   2. Remove the stray <w:lastRenderedPageBreak/> before the "Grading" heading
#   3. Change "half credit" -> "credit" in the late-journal policy paragraph
#   4. Split the "Other notes" sentence so "Thus" is wrapped in gramStart/gramEnd
#      proofing-error markers (as Word's grammar checker would do)
#   5. Update the "Revised:" date field cached text in the footer

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert "Please do not try to link a file..." paragraph after the
#    "Journals will be submitted via a text box on Canvas." paragraph.
# ---------------------------------------------------------------------------
$submissionPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Journals will be submitted via a text box on Canvas.*") {
        $submissionPara = $p
    }
}
if ($submissionPara -ne $null) {
    $submissionPara.Range.InsertParagraphAfter()
    $paras = $d.Paragraphs
    foreach ($p in $paras) {
        if ($p.Range.Text -like "Journals will be submitted via a text box on Canvas.*") {
            $afterPara = $p.Next()
        }
    }
    $afterPara.Range.InsertParagraphAfter()
    $afterPara2 = $afterPara.Next()
    $afterPara2.Range.InsertAfter("Please do not try to link a file by copying the link into the submission box, it will not work.")
}

# ---------------------------------------------------------------------------
# 2. Drop the stray <w:lastRenderedPageBreak/> marker in front of "Grading "
# ---------------------------------------------------------------------------
$gradingPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Grading `r") {
        $gradingPara = $p
    }
}
if ($gradingPara -ne $null) {
    $contentRange = $d.Range($gradingPara.Range.Start, $gradingPara.Range.End - 1)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Grading </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $contentRange.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 3. "I will accept a late journal for half credit up until ..." -> drop "half "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("half credit", $false, $true, $false, $false, $false, $true, 1, $false, "credit", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Wrap "Thus" with gramStart/gramEnd proofErr markers, splitting the run
# ---------------------------------------------------------------------------
$otherNotesPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "When time permits*") {
        $otherNotesPara = $p
    }
}
if ($otherNotesPara -ne $null) {
    $contentRange = $d.Range($otherNotesPara.Range.Start, $otherNotesPara.Range.End - 1)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">When time permits, I will try to answer questions raised in journals. However, I often get behind in grading journals. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Thus</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> any question which you really want answered should be asked via email or Canvas.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $contentRange.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 5. Update the cached "Revised:" DATE field text in the footer
# ---------------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$footer = $sec.Footers.Item(1)
if ($footer.Exists) {
    $footer.Range.Find.Execute("2020-08-21", $false, $true, $false, $false, $false, $true, 1, $false, "2022-08-28", 2) | Out-Null
}
